$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four data rows (9-12) got cyclically shifted down by one record:
# new row9 <- old row12, new row10 <- old row9, new row11 <- old row10, new row12 <- old row11.
# Only columns A,B,D,E,F,G,H,Q,R,Y,AA actually differ between the source records.

# Row 9 (was row 12's record)
$ws.Range("A9").Value = 105343625
$ws.Range("Q9").Value = 579018.9757471241
$ws.Range("R9").Value = 7015651.889117917
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2022-11-02"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2022-11-02"

# Row 10 (was row 9's record)
$ws.Range("A10").Value = 105343678
$ws.Range("B10").Value = 77506
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("Q10").Value = 579367.4720198719
$ws.Range("R10").Value = 7015752.363567609

# Row 11 (was row 10's record)
$ws.Range("A11").Value = 105343692
$ws.Range("B11").Value = 89356
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 5447
$ws.Range("F11").Value = "Vedticka"
$ws.Range("G11").Value = "Fuscoporia viticola"
$ws.Range("H11").Value = "(Schwein.) Murrill"
$ws.Range("Q11").Value = 579369.0961422946
$ws.Range("R11").Value = 7015814.571333516

# Row 12 (was row 11's record)
$ws.Range("A12").Value = 105343677
$ws.Range("Q12").Value = 579313.0396324483
$ws.Range("R12").Value = 7015709.576337469
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = "2022-11-03"
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = "2022-11-03"
